# Apply updated dSF (column F) values for several rows as part of the
# "repull data, push all data, mean calculation" update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = -2
$ws.Range("F14").Value = -2
$ws.Range("F16").Value = -3
$ws.Range("F17").Value = -8
$ws.Range("F22").Value = -5
